$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: set E1 = "ERROR" with same style as D1 (bold, centered, bordered)
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 5).PasteSpecial(-4122)
$ws.Cells.Item(1, 5).Value = "ERROR"

# Copy the column-A style (bold, centered, bordered) onto the newly added rows 56-58
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A56:A58").PasteSpecial(-4122)

# Write data rows 2..58 (A:E) with final values
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 1063.9
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 0.1490276777973729
$ws.Cells.Item(2, 5).Value = 0.007251346738462776
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1085
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 0.2937404414934305
$ws.Cells.Item(3, 5).Value = 0.002862269831848287
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1175
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0.03502302754510135
$ws.Cells.Item(4, 5).Value = 0.0009639365379385692
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 1209.564221739463
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 0.07795795274989076
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1209.564221739463
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 4).Value = 0.01309811992587076
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 1248
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 0.07815507531638952
$ws.Cells.Item(7, 5).Value = 0.001450300366695888
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 1266.460117430137
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0.0008324633664117887
$ws.Cells.Item(8, 5).Value = 0.0002601448020036839
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 1292.473979699731
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = 0.02049772367939864
$ws.Cells.Item(9, 5).Value = 0.0008710852028429702
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1304.585912295229
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 0.02139510558063962
$ws.Cells.Item(10, 5).Value = 0.0009290339588017206
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 1317.3
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.004773763494720345
$ws.Cells.Item(11, 5).Value = 0.000429102336604076
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 1339.018418868032
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = 0.007785429305680182
$ws.Cells.Item(12, 5).Value = 0.0005245502014465371
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 1358.3
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = 0.4491463597744076
$ws.Cells.Item(13, 5).Value = 0.003666727657259063
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 1456.533865675461
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(14, 4).Value = 0.08340446253726887
$ws.Cells.Item(14, 5).Value = 0.00583831237760882
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 1471.462544459845
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 4).Value = 0.06808118586513147
$ws.Cells.Item(15, 5).Value = 0.003275876521375345
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 1494.44871075181
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 0.1900417373909443
$ws.Cells.Item(16, 5).Value = 0.01181851600689953
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 1515.596556909176
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0.2269019553533708
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 1515.596556909176
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 0.03534896350903715
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 1538.408632604229
$ws.Cells.Item(19, 3).Value = 4
$ws.Cells.Item(19, 4).Value = 0.04823577340117313
$ws.Cells.Item(19, 5).Value = 0.004823577340117314
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 1574
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 0.07289689843887331
$ws.Cells.Item(20, 5).Value = 0.009183861220645457
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 1597
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = 0.4214536699028364
$ws.Cells.Item(21, 5).Value = 0.01760131336021126
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 1613.149084845743
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0.003351768912339795
$ws.Cells.Item(22, 5).Value = 0.0003511376955784547
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 1625.344244785309
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 0.004628586547637729
$ws.Cells.Item(23, 5).Value = 0.0006612266496625326
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 1649.346715797556
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0.0270023364094122
$ws.Cells.Item(24, 5).Value = 0.0009915178608107295
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 1732.134406941665
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0.003175346061097805
$ws.Cells.Item(25, 5).Value = 0.0004885147786304314
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 1747.835619913894
$ws.Cells.Item(26, 3).Value = 4
$ws.Cells.Item(26, 4).Value = 0.05870067363643368
$ws.Cells.Item(26, 5).Value = 0.00571955281585764
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 1772.725731049579
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 0.006562677749158311
$ws.Cells.Item(27, 5).Value = 0.0004732700299873783
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 1785.161835361124
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = 0.03823531669500906
$ws.Cells.Item(28, 5).Value = 0.001203763077390509
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 1795.477475116088
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = 0.01013878774325021
$ws.Cells.Item(29, 5).Value = 0.0007310097745917075
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 1805.922029018604
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = 0.0507632814884952
$ws.Cells.Item(30, 5).Value = 0.004587043507996554
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 1818
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0.08069046492397715
$ws.Cells.Item(31, 5).Value = 0.002006116531259101
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 1840
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 0.06460444685554217
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 1840
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 0.01312801055132022
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 1841.298544231364
$ws.Cells.Item(34, 3).Value = 2
$ws.Cells.Item(34, 4).Value = 0.002547702085345394
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 1841.298544231364
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 0.002652847806461739
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 1851.005074483267
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 0.005867300457096876
$ws.Cells.Item(36, 5).Value = 0.0005333909506451705
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 1865.530752678807
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(37, 4).Value = 0.05153847163427385
$ws.Cells.Item(37, 5).Value = 0.001688721187764487
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 1876
$ws.Cells.Item(38, 3).Value = 2
$ws.Cells.Item(38, 4).Value = 0.04793841377196112
$ws.Cells.Item(38, 5).Value = 0.001727216811634221
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 1885.117065551615
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(39, 4).Value = 0.1510466205143166
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 1885.117065551615
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 0.003039328872354414
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 1896.361077451411
$ws.Cells.Item(41, 3).Value = 2
$ws.Cells.Item(41, 4).Value = 0.007999904990412195
$ws.Cells.Item(41, 5).Value = 0.000557370429659866
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 1911.792281651251
$ws.Cells.Item(42, 3).Value = 2
$ws.Cells.Item(42, 4).Value = 0.01694965968195618
$ws.Cells.Item(42, 5).Value = 0.0009544468558771441
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 1924.67520477692
$ws.Cells.Item(43, 3).Value = 2
$ws.Cells.Item(43, 4).Value = 0.01699344876080397
$ws.Cells.Item(43, 5).Value = 0.004760810547773875
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 1933.550109340654
$ws.Cells.Item(44, 3).Value = 2
$ws.Cells.Item(44, 4).Value = 0.01755077824842062
$ws.Cells.Item(44, 5).Value = 0.004047627497745399
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 1944.234562301843
$ws.Cells.Item(45, 3).Value = 2
$ws.Cells.Item(45, 4).Value = 0.01052252851291888
$ws.Cells.Item(45, 5).Value = 0.001698256184047034
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 1959.440632592807
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 0.01401285418419295
$ws.Cells.Item(46, 5).Value = 0.001137081246450025
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 1970.376107067389
$ws.Cells.Item(47, 3).Value = 2
$ws.Cells.Item(47, 4).Value = 0.01625312779948391
$ws.Cells.Item(47, 5).Value = 0.001175329489632101
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 1986.766402625133
$ws.Cells.Item(48, 3).Value = 2
$ws.Cells.Item(48, 4).Value = 0.01049407157396844
$ws.Cells.Item(48, 5).Value = 0.0007760888945378589
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 1995.02814271466
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = 0.00358591301441812
$ws.Cells.Item(49, 5).Value = 0.0002098721842821674
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 2014.856993695083
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(50, 4).Value = 0.001969364201798255
$ws.Cells.Item(50, 5).Value = 0.0001188354704612277
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 2048.833064430492
$ws.Cells.Item(51, 3).Value = 2
$ws.Cells.Item(51, 4).Value = 0.03176164738981578
$ws.Cells.Item(51, 5).Value = 0.001065596395112867
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 2073.838281975398
$ws.Cells.Item(52, 3).Value = 2
$ws.Cells.Item(52, 4).Value = 0.04098796786726889
$ws.Cells.Item(52, 5).Value = 0.001247309250398715
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 2092.576632909581
$ws.Cells.Item(53, 3).Value = 2
$ws.Cells.Item(53, 4).Value = 0.02568909853921173
$ws.Cells.Item(53, 5).Value = 0.00101083291402597
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(54, 2).Value = 2113.24
$ws.Cells.Item(54, 3).Value = 2
$ws.Cells.Item(54, 4).Value = 0.5303015044663604
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(55, 2).Value = 2113.24
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 0.01248143625838344
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(56, 2).Value = 2133.721396755173
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 0.005608549416656394
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(57, 2).Value = 2133.721396755173
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(57, 4).Value = 0.2481396455136796
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = 2148.119842376009
$ws.Cells.Item(58, 3).Value = 2
$ws.Cells.Item(58, 4).Value = 0.01106723146435628
$ws.Cells.Item(58, 5).Value = 0.0007449098101009035

Write-Output "done"
